# Actualización automática 2025-12-08 14:30:06
#
# Updates sales figures for ILLER LOPEZ ROBERTO FERNANDO / SARMIENTO SARMIENTO
# SANDRA EULALIA on both sheets of the workbook.

$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" (row 19 = client data, row 21 = "x de 19" counters)
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

$wsGrupo.Range("E19").Value = 643.6
$wsGrupo.Range("M19").Value = 2319.99

$wsGrupo.Range("E21").Value = "2 de 19"
$wsGrupo.Range("M21").Value = "1 de 19"

# --- Sheet "VENTA MENSUAL" (row 19 = client data, row 21 = column totals)
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

$wsMensual.Range("F19").Value = 2963.59
$wsMensual.Range("F21").Value = 2465.97
